$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A200").Value = "TestVal"
$v2 = $ws.Range("A200").Value2
Write-Output ("A200:" + $v2)
